# Update attendee/ticket numbers ("想去人数" column F) on the "展览" and
# "全部类型" worksheets to reflect the latest scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 4-9, column F) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 52
$ws1.Range("F5").Value = 5010
$ws1.Range("F6").Value = 171
$ws1.Range("F7").Value = 85
$ws1.Range("F8").Value = 297
$ws1.Range("F9").Value = 49

# --- Sheet "全部类型" (rows 8-14, column F, skipping row 12) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 52
$ws4.Range("F9").Value = 5010
$ws4.Range("F10").Value = 171
$ws4.Range("F11").Value = 85
$ws4.Range("F13").Value = 297
$ws4.Range("F14").Value = 49
